# Commit: "update to manual status column;"
#
# The I column ("manualStatus") previously held raw numeric flag codes
# (128, or 32128 for row 26). It now holds the equivalent bracketed,
# comma-separated text label instead, e.g. 128 -> "[128]" and
# 32128 -> "[32,128]".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = "[128]"
}
$ws.Cells.Item(26, 9).Value = "[32,128]"
$ws.Cells.Item(27, 9).Value = "[128]"

[void]$ws.Range("I27").Select()
